$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data values in A2:B5 (new cluster analysis results)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 1836

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1202

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 1045

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 162

# Remove row 6 entirely, shrinking the used range to A1:B5
$ws.Range("A6:B6").Delete()
